# Negate all numeric values in column E ("Block") across the data rows.
# This mirrors the diff where every <c r="E.." t="n"><v>N</v></c> had its
# sign flipped (N -> -N).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = 157

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 5)
    $val = $cell.Value2
    if ($val -is [double] -or $val -is [int]) {
        $cell.Value2 = (-1) * $val
    }
}
